$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bug列表")

# Insert a new row at row 62, shifting rows 62:117 down to 63:118
$ws.Rows("62:62").Insert()

$ws.Cells.Item(62, 4).Value = "回购审批，单独的回购审批页面"
$ws.Cells.Item(62, 5).Value = "易用性"
$ws.Cells.Item(62, 6).Value = "Tina"
